$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value, derived from the authoritative diff.
# NumberFormat is forced to Text ("@") before assignment so that numeric-
# looking strings (e.g. '1.00', '0.576', '4.109.00') are preserved verbatim
# as text instead of being auto-coerced into numeric values by Excel.
$updates = [ordered]@{
    'D2' = '66.979.34'
    'E2' = '  +5.13%  '
    'D3' = '3.506.99'
    'E3' = '  +2.85%  '
    'E4' = '  -0.01%  '
    'D5' = '594.99'
    'E5' = '  +4.51%  '
    'D6' = '169.43'
    'E6' = '  +7.48%  '
    'E7' = '  -0.04%  '
    'D8' = '3.506.30'
    'E8' = '  +2.71%  '
    'D9' = '0.576'
    'E9' = '  +1.54%  '
    'E10' = '  +0.69%  '
    'E11' = '  +5.67%  '
    'E12' = '  +4.15%  '
    'D13' = '4.109.00'
    'E14' = '  +0.12%  '
    'D15' = '28.26'
    'E15' = '  +4.59%  '
    'E16' = '  +4.61%  '
    'D17' = '66.905.65'
    'E17' = '  +4.90%  '
    'D18' = '3.511.58'
    'E18' = '  +2.60%  '
    'E19' = '  +3.80%  '
    'E20' = '  +3.34%  '
    'D21' = '395.13'
    'E21' = '  +2.50%  '
    'E22' = '  +2.41%  '
    'D23' = '73.22'
    'E23' = '  +2.81%  '
    'E24' = '  +12.20%  '
    'D26' = '0.531'
    'E26' = '  +3.10%  '
    'D27' = '10.07'
    'E27' = '  +4.06%  '
    'E28' = '  +2.55%  '
    'E29' = '  +0.13%  '
    'E30' = '  +4.66%  '
    'E31' = '  +6.05%  '
    'D32' = '2.07'
    'E32' = '  +4.54%  '
    'D33' = '23.57'
    'E33' = '  +3.05%  '
    'D34' = '7.47'
    'E34' = '  +7.45%  '
    'D35' = '1.00'
    'E35' = '  +0.08%  '
    'E36' = '  +6.24%  '
    'D37' = '162.57'
    'E37' = '  +1.28%  '
    'D38' = '0.902'
    'E38' = '  +6.77%  '
    'E39' = '  +6.75%  '
    'D40' = '0.0754'
    'D41' = '4.68'
    'E41' = '  +7.57%  '
    'B42' = 'Maker'
    'C42' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D42' = '2.842.35'
    'E42' = '  +1.86%  '
    'B43' = 'RenderToken'
    'C43' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D43' = '6.68'
    'E43' = '  +5.02%  '
    'E44' = '  +1.85%  '
    'D45' = '43.51'
    'E45' = '  +1.19%  '
    'D46' = '26.59'
    'E46' = '  +2.34%  '
    'E47' = '  +4.55%  '
    'D48' = '2.56'
    'E48' = '  +7.84%  '
    'D49' = '349.06'
    'E50' = '  +5.06%  '
    'D51' = '33.76'
    'E51' = '  +12.37%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
}
